# remove index for title
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading "N. " numbering prefix from every book title in column A
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $title = $cell.Value2
    $newTitle = $title -replace '^\d+\.\s', ''
    $cell.Value = $newTitle
}

# Helper to write a value as a real text string (not auto-converted to a number)
function Set-TextValue($range, $text) {
    $range.Formula = "=""" + $text + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# A handful of rating counts were also refreshed alongside the title edits
Set-TextValue $ws.Range("D4") "7 "
Set-TextValue $ws.Range("D10") "28,309 "
Set-TextValue $ws.Range("D17") "38,345 "
Set-TextValue $ws.Range("D20") "183,328 "
Set-TextValue $ws.Range("D31") "588 "

$excel.CutCopyMode = $false
